$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Hunk 1: row 10 (Item ID 1959)
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").Value = $null

# Hunk 2: row 33 (Item ID 5512)
$ws.Range("H33").Value = 959.96875
$ws.Range("I33").Value = 817.5
$ws.Range("K33").Value = 817.5
$ws.Range("M33").Value = -588.5

# Hunk 3: row 116 (Item ID 27778)
$ws.Range("H116").Value = 3048.2
$ws.Range("I116").Value = 3074.1538
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 3074.1538
$ws.Range("L116").Value = 3000
$ws.Range("M116").Value = 367.8462
$ws.Range("N116").Value = -9884

# Hunk 4: row 138 (Item ID 44169)
$ws.Range("H138").Value = 2326.7673
$ws.Range("J138").Value = 2069.1724
$ws.Range("L138").Value = 6207.5172
$ws.Range("N138").Value = -16487.5172

# Hunk 5: row 140 (Item ID 42459)
$ws.Range("H140").Value = 40780
$ws.Range("J140").Value = 40780
$ws.Range("L140").Value = 40780
$ws.Range("N140").Value = -51140

$ws = $wb.Worksheets.Item("ARM")
# Hunk 6: row 28 (Item ID 19534)
$ws.Range("H28").Value = 3865.25
$ws.Range("I28").Value = 3865.25
$ws.Range("K28").Value = 3865.25
$ws.Range("M28").Value = -3673.25

# Hunk 7: row 99 (Item ID 19534)
$ws.Range("H99").Value = 3865.25
$ws.Range("I99").Value = 3865.25
$ws.Range("K99").Value = 3865.25
$ws.Range("M99").Value = -870.25

# Hunk 8: row 108 (Item ID 27084)
$ws.Range("H108").Value = 31000
$ws.Range("J108").Value = 31000
$ws.Range("L108").Value = 31000
$ws.Range("N108").Value = -38680

# Hunk 9: row 122 (Item ID 36168)
$ws.Range("H122").Value = 2460.182
$ws.Range("I122").Value = 2069.6924
$ws.Range("J122").Value = 3024.2222
$ws.Range("K122").Value = 6209.0772
$ws.Range("L122").Value = 9072.6666
$ws.Range("M122").Value = -3759.0772
$ws.Range("N122").Value = -13972.6666

# Hunk 10: row 123 (Item ID 34107)
$ws.Range("H123").Value = 39429
$ws.Range("J123").Value = 39429
$ws.Range("L123").Value = 39429
$ws.Range("N123").Value = -49229

$ws = $wb.Worksheets.Item("BSM")
# Hunk 11: row 50 (Item ID 27159)
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").Value = $null

# Hunk 12: row 97 (Item ID 19518)
$ws.Range("H97").Value = 41157
$ws.Range("I97").Value = 20000
$ws.Range("J97").Value = 51735.5
$ws.Range("K97").Value = 20000
$ws.Range("L97").Value = 51735.5
$ws.Range("M97").Value = -19009
$ws.Range("N97").Value = -53717.5

# Hunk 13: row 109 (Item ID 27096)
$ws.Range("H109").Value = 37272.727
$ws.Range("J109").Value = 37272.727
$ws.Range("L109").Value = 37272.727
$ws.Range("N109").Value = -40046.727

# Hunk 14: row 115 (Item ID 27118)
$ws.Range("H115").Value = 69342
$ws.Range("J115").Value = 69342
$ws.Range("L115").Value = 69342
$ws.Range("N115").Value = -72476

# Hunk 15: row 118 (Item ID 27137)
$ws.Range("H118").Value = 41785.332
$ws.Range("J118").Value = 41785.332
$ws.Range("L118").Value = 41785.332
$ws.Range("N118").Value = -45099.332

# Hunk 16: row 135 (Item ID 41992)
$ws.Range("H135").Value = 69726
$ws.Range("I135").Value = 80000
$ws.Range("J135").Value = 68992.14
$ws.Range("K135").Value = 80000
$ws.Range("L135").Value = 68992.14
$ws.Range("M135").Value = -74930
$ws.Range("N135").Value = -79132.14

$ws = $wb.Worksheets.Item("CRP")
# Hunk 17: row 10 (Item ID 1997)
$ws.Range("H10").Value = 18333.334

# Hunk 18: row 106 (Item ID 18661)
$ws.Range("H106").Value = 33500
$ws.Range("J106").Value = 33500
$ws.Range("L106").Value = 33500
$ws.Range("N106").Value = -36024

# Hunk 19: row 112 (Item ID 25796)
$ws.Range("H112").Value = 39000
$ws.Range("J112").Value = 39000
$ws.Range("L112").Value = 39000
$ws.Range("N112").Value = -41954

# Hunk 20: row 114 (Item ID 27112)
$ws.Range("H114").Value = 57228
$ws.Range("J114").Value = 57228
$ws.Range("L114").Value = 57228
$ws.Range("N114").Value = -65906

$ws = $wb.Worksheets.Item("CUL")
# Hunk 21: row 17 (Item ID 4640)
$ws.Range("H17").Value = 482.22223
$ws.Range("I17").Value = 486.25
$ws.Range("K17").Value = 1458.75
$ws.Range("M17").Value = -1289.75

# Hunk 22: row 117 (Item ID 27870)
$ws.Range("H117").Value = 459
$ws.Range("I117").Value = 459
$ws.Range("K117").Value = 1377
$ws.Range("M117").Value = 2065

$ws = $wb.Worksheets.Item("GSM")
# Hunk 23: row 39 (Item ID 18264)
$ws.Range("H39").Value = 23000
$ws.Range("J39").Value = 23000
$ws.Range("L39").Value = 23000
$ws.Range("N39").Value = -24064

# Hunk 24: row 62 (Item ID 11983)
$ws.Range("H62").Value = 30000
$ws.Range("J62").Value = 30000
$ws.Range("L62").Value = 30000
$ws.Range("N62").Value = -31372

# Hunk 25: row 65 (Item ID 11983)
$ws.Range("H65").Value = 30000
$ws.Range("J65").Value = 30000
$ws.Range("L65").Value = 90000
$ws.Range("N65").Value = -96864

# Hunk 26: row 92 (Item ID 18094)
$ws.Range("H92").Value = 7000
$ws.Range("J92").Value = 7000
$ws.Range("L92").Value = 7000
$ws.Range("N92").Value = -10744

# Hunk 27: row 99 (Item ID 19532)
$ws.Range("H99").Value = 24111.766
$ws.Range("I99").Value = 30000
$ws.Range("J99").Value = 23743.75
$ws.Range("K99").Value = 30000
$ws.Range("L99").Value = 23743.75
$ws.Range("M99").Value = -27754
$ws.Range("N99").Value = -28235.75

# Hunk 28: row 122 (Item ID 36182)
$ws.Range("H122").Value = 5305.517
$ws.Range("I122").Value = 1833.3334
$ws.Range("J122").Value = 5706.154
$ws.Range("K122").Value = 5500.0002
$ws.Range("L122").Value = 17118.462
$ws.Range("M122").Value = -3050.0002
$ws.Range("N122").Value = -22018.462

# Hunk 29: row 130 (Item ID 34692)
$ws.Range("H130").Value = 58800
$ws.Range("J130").Value = 58800
$ws.Range("L130").Value = 58800
$ws.Range("N130").Value = -68840

$ws = $wb.Worksheets.Item("LTW")
# Hunk 30: row 14 (Item ID 3771)
$ws.Range("H14").Value = 30510.9
$ws.Range("I14").Value = 14552
$ws.Range("J14").Value = 34500.625
$ws.Range("K14").Value = 14552
$ws.Range("L14").Value = 34500.625
$ws.Range("M14").Value = -14380
$ws.Range("N14").Value = -34844.625

# Hunk 31: row 122 (Item ID 36247)
$ws.Range("H122").Value = 2653
$ws.Range("I122").Value = 2479.5
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 7438.5
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -4988.5
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("WVR")
# Hunk 32: row 46 (Item ID 42037)
$ws.Range("H46").Value = 47628.43
$ws.Range("J46").Value = 47628.43
$ws.Range("L46").Value = 47628.43
$ws.Range("N46").Value = -48090.43

# Hunk 33: row 97 (Item ID 18220)
$ws.Range("H97").Value = 43690.668
$ws.Range("J97").Value = 43690.668
$ws.Range("L97").Value = 43690.668
$ws.Range("N97").Value = -45672.668

# Hunk 34: row 103 (Item ID 18548)
$ws.Range("H103").Value = 10602
$ws.Range("J103").Value = 10602
$ws.Range("L103").Value = 10602
$ws.Range("N103").Value = -12946

# Hunk 35: row 130 (Item ID 34705)
$ws.Range("H130").Value = 53595.668
$ws.Range("J130").Value = 53595.668
$ws.Range("L130").Value = 53595.668
$ws.Range("N130").Value = -63635.668

# Hunk 36: row 134 (Item ID 42037)
$ws.Range("H134").Value = 47628.43
$ws.Range("J134").Value = 47628.43
$ws.Range("L134").Value = 142885.29
$ws.Range("N134").Value = -147955.29

